$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark from its current location (the paragraph
#    holding the inline picture, right before the trailing empty paragraphs).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Append three new "Normal (Web)" paragraphs of text plus two new empty
#    "Normal (Web)" paragraphs (mirroring the very last paragraph's
#    formatting) at the end of the document body, re-homing the _GoBack
#    bookmark into the new final paragraph.
#
#    InsertXML on a range collapsed at the very end of the document always
#    folds the *last* supplied <w:p> into the document's existing final
#    paragraph mark (keeping that paragraph's own pPr, only adopting any
#    runs/bookmarks we give it). So we give it a totally empty trailing
#    <w:p/> placeholder to fold into unchanged, and everything before that
#    becomes its own brand-new, independently-formatted paragraph.
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$para1 = "<w:p $ns>" +
  "<w:pPr>" +
    "<w:pStyle w:val='5'/>" +
    "<w:keepNext w:val='0'/>" +
    "<w:keepLines w:val='0'/>" +
    "<w:widowControl/>" +
    "<w:suppressLineNumbers w:val='0'/>" +
    "<w:spacing w:before='0' w:beforeAutospacing='0' w:after='0' w:afterAutospacing='0' w:line='380' w:lineRule='atLeast'/>" +
    "<w:ind w:left='0' w:right='0'/>" +
    "<w:jc w:val='left'/>" +
    "<w:rPr>" +
      "<w:rFonts w:ascii='helvetica neue' w:hAnsi='helvetica neue' w:eastAsia='helvetica neue' w:cs='helvetica neue'/>" +
      "<w:kern w:val='0'/>" +
      "<w:sz w:val='26'/>" +
      "<w:szCs w:val='26'/>" +
      "<w:lang w:val='en-US' w:eastAsia='zh-CN' w:bidi='ar'/>" +
    "</w:rPr>" +
  "</w:pPr>" +
  "<w:r>" +
    "<w:rPr>" +
      "<w:rFonts w:ascii='helvetica neue' w:hAnsi='helvetica neue' w:eastAsia='helvetica neue' w:cs='helvetica neue'/>" +
      "<w:kern w:val='0'/>" +
      "<w:sz w:val='26'/>" +
      "<w:szCs w:val='26'/>" +
      "<w:lang w:val='en-US' w:eastAsia='zh-CN' w:bidi='ar'/>" +
    "</w:rPr>" +
    "<w:t>想个名字呢？</w:t>" +
  "</w:r>" +
"</w:p>"

$para2 = "<w:p $ns>" +
  "<w:pPr>" +
    "<w:pStyle w:val='5'/>" +
    "<w:keepNext w:val='0'/>" +
    "<w:keepLines w:val='0'/>" +
    "<w:widowControl/>" +
    "<w:suppressLineNumbers w:val='0'/>" +
    "<w:spacing w:before='0' w:beforeAutospacing='0' w:after='0' w:afterAutospacing='0' w:line='380' w:lineRule='atLeast'/>" +
    "<w:ind w:left='0' w:right='0' w:firstLine='420' w:firstLineChars='0'/>" +
    "<w:jc w:val='left'/>" +
  "</w:pPr>" +
  "<w:r>" +
    "<w:rPr>" +
      "<w:rFonts w:ascii='helvetica neue' w:hAnsi='helvetica neue' w:eastAsia='helvetica neue' w:cs='helvetica neue'/>" +
      "<w:kern w:val='0'/>" +
      "<w:sz w:val='26'/>" +
      "<w:szCs w:val='26'/>" +
      "<w:lang w:val='en-US' w:eastAsia='zh-CN' w:bidi='ar'/>" +
    "</w:rPr>" +
    "<w:t>极限编程中所说的“隐喻”</w:t>" +
  "</w:r>" +
"</w:p>"

$para3 = "<w:p $ns>" +
  "<w:pPr>" +
    "<w:pStyle w:val='5'/>" +
    "<w:keepNext w:val='0'/>" +
    "<w:keepLines w:val='0'/>" +
    "<w:widowControl/>" +
    "<w:suppressLineNumbers w:val='0'/>" +
    "<w:spacing w:before='0' w:beforeAutospacing='0' w:after='0' w:afterAutospacing='0' w:line='380' w:lineRule='atLeast'/>" +
    "<w:ind w:left='0' w:right='0' w:firstLine='420' w:firstLineChars='0'/>" +
    "<w:jc w:val='left'/>" +
  "</w:pPr>" +
  "<w:r>" +
    "<w:rPr>" +
      "<w:rFonts w:ascii='helvetica neue' w:hAnsi='helvetica neue' w:eastAsia='helvetica neue' w:cs='helvetica neue'/>" +
      "<w:kern w:val='0'/>" +
      "<w:sz w:val='26'/>" +
      "<w:szCs w:val='26'/>" +
      "<w:lang w:val='en-US' w:eastAsia='zh-CN' w:bidi='ar'/>" +
    "</w:rPr>" +
    "<w:t>找到合适的名字，理解各个东西是怎么彼此关联的</w:t>" +
  "</w:r>" +
"</w:p>"

$para4 = "<w:p $ns>" +
  "<w:pPr>" +
    "<w:pStyle w:val='5'/>" +
    "<w:keepNext w:val='0'/>" +
    "<w:keepLines w:val='0'/>" +
    "<w:widowControl/>" +
    "<w:suppressLineNumbers w:val='0'/>" +
    "<w:pBdr>" +
      "<w:top w:val='none' w:color='auto' w:sz='0' w:space='0'/>" +
      "<w:left w:val='none' w:color='auto' w:sz='0' w:space='0'/>" +
      "<w:bottom w:val='none' w:color='auto' w:sz='0' w:space='0'/>" +
      "<w:right w:val='none' w:color='auto' w:sz='0' w:space='0'/>" +
    "</w:pBdr>" +
    "<w:spacing w:before='0' w:beforeAutospacing='0' w:after='300' w:afterAutospacing='0'/>" +
    "<w:ind w:right='0'/>" +
    "<w:jc w:val='left'/>" +
    "<w:rPr>" +
      "<w:rFonts w:hint='default' w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
      "<w:i w:val='0'/>" +
      "<w:caps w:val='0'/>" +
      "<w:color w:val='333333'/>" +
      "<w:spacing w:val='0'/>" +
      "<w:sz w:val='32'/>" +
      "<w:szCs w:val='32'/>" +
    "</w:rPr>" +
  "</w:pPr>" +
"</w:p>"

$para5 = "<w:p $ns>" +
  "<w:pPr>" +
    "<w:pStyle w:val='5'/>" +
    "<w:keepNext w:val='0'/>" +
    "<w:keepLines w:val='0'/>" +
    "<w:widowControl/>" +
    "<w:suppressLineNumbers w:val='0'/>" +
    "<w:pBdr>" +
      "<w:top w:val='none' w:color='auto' w:sz='0' w:space='0'/>" +
      "<w:left w:val='none' w:color='auto' w:sz='0' w:space='0'/>" +
      "<w:bottom w:val='none' w:color='auto' w:sz='0' w:space='0'/>" +
      "<w:right w:val='none' w:color='auto' w:sz='0' w:space='0'/>" +
    "</w:pBdr>" +
    "<w:spacing w:before='0' w:beforeAutospacing='0' w:after='300' w:afterAutospacing='0'/>" +
    "<w:ind w:right='0'/>" +
    "<w:jc w:val='left'/>" +
    "<w:rPr>" +
      "<w:rFonts w:hint='default' w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
      "<w:i w:val='0'/>" +
      "<w:caps w:val='0'/>" +
      "<w:color w:val='333333'/>" +
      "<w:spacing w:val='0'/>" +
      "<w:sz w:val='32'/>" +
      "<w:szCs w:val='32'/>" +
    "</w:rPr>" +
  "</w:pPr>" +
  "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
  "<w:bookmarkEnd w:id='0'/>" +
"</w:p>"

# Trailing empty placeholder paragraph: InsertXML folds the last <w:p> it is
# given into the pre-existing final paragraph mark (keeping that paragraph's
# own pPr). Supplying it empty means that pre-existing last paragraph is
# left completely untouched, as required by the diff.
$placeholder = "<w:p $ns></w:p>"

$xml = $para1 + $para2 + $para3 + $para4 + $para5 + $placeholder

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML($xml)
